$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the existing Hyperlink cell style (xf index 2, "Hyperlink" named style)
# from one of the current hyperlink cells before we touch anything, so we can
# re-apply it later without Excel fabricating a brand-new style entry.
$hlStyle = $ws.Range("F2").Style

# The "Language Used" column (C) is being removed entirely; drop the existing
# hyperlinks first since this engine does not auto-shift hyperlink ranges when
# columns/rows are deleted - we'll re-create them at their new locations below.
$ws.Hyperlinks.Delete()

# Remove column C ("Language Used") - this shifts D,E,F,G left to C,D,E,F.
$ws.Columns("C").Delete()

# Remove the now-stray row 5 (previously just held a trailing Serial No of 4
# with no other data).
$ws.Rows("5").Delete()

# Re-create the hyperlinks on their new (shifted-left) cells.
$ws.Hyperlinks.Add($ws.Range("E2"), "https://github.com/abhisekjha/steganography") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.abhisekjha.com.np/steganography/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("E3"), "https://github.com/abhisekjha/pqc_aes_multipath") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F3"), "https://www.abhisekjha.com.np/pqc_aes_multipath") | Out-Null
$ws.Hyperlinks.Add($ws.Range("E4"), "https://github.com/abhisekjha/luxury-car-calculator") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F4"), "https://www.abhisekjha.com.np/luxury-car-calculator") | Out-Null

# Hyperlinks.Add() applies its own freshly-built style; restore the original
# shared Hyperlink style so we don't leave an extra duplicate style behind.
$ws.Range("E2").Style = $hlStyle
$ws.Range("F2").Style = $hlStyle
$ws.Range("E3").Style = $hlStyle
$ws.Range("F3").Style = $hlStyle
$ws.Range("E4").Style = $hlStyle
$ws.Range("F4").Style = $hlStyle

# Update the selected cell to match the author's final cursor position.
$ws.Range("A4").Select()
